# Edit script: (1) re-apply the built-in table style to the three summary
# tables on slides 14-16, (2) swap the "Integral" / "Office Theme" colour
# themes between the slide master and the notes master (theme1.xml <->
# theme2.xml content swap).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Tables: switch tableStyleId from {297E6AD9-...} to {75343CEB-...}
# ---------------------------------------------------------------------
$tableSlideIndexes = @(14, 15, 16)
foreach ($idx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($idx)
    $tblShape = $slide.Shapes.Item(1)
    if ($tblShape.HasTable) {
        $tblShape.Table.ApplyStyle("{75343CEB-610D-43A3-B93B-F74B8D59EEB8}")
    }
}

# ---------------------------------------------------------------------
# 2. Theme swap: Master currently uses the "Integral" / "Red Violet"
#    colours, NotesMaster currently uses "Office Theme" / "Office"
#    colours. Swap them so the slide master becomes the default Office
#    colours and the notes master becomes the Integral colours.
# ---------------------------------------------------------------------

# RGB() equivalents (VBA long = B*65536 + G*256 + R) for each theme,
# in ColorScheme.Item() order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)
$integralColors = @(0, 16777215, 5326149, 14473688, 9514467, 13381832, 14460494, 15168839, 14774665, 7555029, 2465643, 9211020)

$master = $p.Slides.Item(1).Master
$masterScheme = $master.ColorScheme
for ($i = 1; $i -le 12; $i++) {
    $masterScheme.Item($i).RGB = $officeColors[$i - 1]
}
$masterScheme.Name = "Office"

$notesMaster = $p.NotesMaster
$notesScheme = $notesMaster.ColorScheme
for ($i = 1; $i -le 12; $i++) {
    $notesScheme.Item($i).RGB = $integralColors[$i - 1]
}
$notesScheme.Name = "Red Violet"
